# Update "想去人数" (interest count, column F) figures on the "展览" and
# "全部类型" sheets to reflect the latest scrape, per commit:
#   "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1798
$ws.Range("F4").Value = 462
$ws.Range("F7").Value = 634
$ws.Range("F8").Value = 344
$ws.Range("F9").Value = 1751
$ws.Range("F10").Value = 374
$ws.Range("F11").Value = 1433
$ws.Range("F12").Value = 817
$ws.Range("F13").Value = 342
$ws.Range("F14").Value = 688
$ws.Range("F15").Value = 12854
$ws.Range("F16").Value = 12841
$ws.Range("F20").Value = 522
$ws.Range("F22").Value = 579
$ws.Range("F23").Value = 2019
$ws.Range("F27").Value = 74
$ws.Range("F29").Value = 682

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 1798
$ws.Range("F6").Value = 462
$ws.Range("F11").Value = 634
$ws.Range("F13").Value = 344
$ws.Range("F14").Value = 1751
$ws.Range("F15").Value = 374
$ws.Range("F16").Value = 1433
$ws.Range("F17").Value = 817
$ws.Range("F18").Value = 342
$ws.Range("F20").Value = 688
$ws.Range("F21").Value = 12854
$ws.Range("F22").Value = 12841
$ws.Range("F26").Value = 522
$ws.Range("F28").Value = 579
$ws.Range("F31").Value = 2019
$ws.Range("F37").Value = 74
$ws.Range("F39").Value = 682
